# validation for user inputs & exceptions
#
# Adds two small reference/demo tables below the existing data on the
# "position" sheet:
#   - rows 75-78: a copy of the First Name/Last name/age/gender block
#     together with a "key1..key4" lookup column (used later for data
#     validation / lookups against user input)
#   - rows 81-85: a "hours per subject per student" table with a bold
#     header row and an AVERAGE() formula summarising each subject

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1 (rows 75-78): First Name / Last name / age / gender, plus the
# key1..key4 helper column in E/F.
# ---------------------------------------------------------------------
$ws.Range("B75").Value = "First Name"
$ws.Range("C75").Value = "Roddy"
$ws.Range("E75").Value = "key1"
$ws.Range("F75").Value = 1

$ws.Range("B76").Value = "Last name"
$ws.Range("C76").Value = "Wiliams"
$ws.Range("E76").Value = "key2"
$ws.Range("F76").Value = 2

$ws.Range("B77").Value = "age"
$ws.Range("C77").Value = 34
$ws.Range("E77").Value = "key3"
$ws.Range("F77").Value = 3

$ws.Range("B78").Value = "gender"
$ws.Range("C78").Value = "male"
$ws.Range("E78").Value = "key4"
$ws.Range("F78").Value = 4

# The "key2"/"key3"/"key4" labels carry a (near-invisible) run-level font
# on the trailing digit, a leftover of how they were originally typed.
# Re-apply the same Arial/10pt font explicitly so the cells get their own
# distinct style from the plain text cells around them.
foreach ($addr in @("E76", "E77", "E78")) {
    $chars = $ws.Range($addr).Characters(4, 1)
    $chars.Font.Name = "Arial"
    $chars.Font.Size = 10
}

# ---------------------------------------------------------------------
# Block 2 (rows 81-85): hours per subject per student, with a bold
# header row and an AVERAGE() summary column.
# ---------------------------------------------------------------------
$ws.Range("C81").Value = "name"
$ws.Range("D81").Value = "hours"
$ws.Range("E81").Value = "Student1"
$ws.Range("F81").Value = "Student2"
$ws.Range("G81").Value = "Student3"
$ws.Range("H81").Value = "sum"
$ws.Range("C81:H81").Font.Bold = $true
$ws.Range("C81:H81").Font.Name = "Arial"
$ws.Range("C81:H81").Font.Size = 10

$ws.Range("C82").Value = "Science"
$ws.Range("D82").Value = 34
$ws.Range("E82").Value = 12
$ws.Range("F82").Value = 5
$ws.Range("G82").Value = 10
$ws.Range("H82").Formula = "=AVERAGE(E82:G82)"

$ws.Range("C83").Value = "ICT"
$ws.Range("D83").Value = 23
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = 9
$ws.Range("G83").Value = 9
$ws.Range("H83").Formula = "=AVERAGE(E83:G83)"

$ws.Range("C84").Value = "History"
$ws.Range("D84").Value = 35
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 11
$ws.Range("G84").Value = 12
$ws.Range("H84").Formula = "=AVERAGE(E84:G84)"

$ws.Range("C85").Value = "Geography"
$ws.Range("D85").Value = 43
$ws.Range("E85").Value = 6
$ws.Range("F85").Value = 5
$ws.Range("G85").Value = 4
$ws.Range("H85").Formula = "=AVERAGE(E85:G85)"

# The subject labels on rows 83-85 wrap within their column (row 82's
# label does not).
$ws.Range("C83:C85").WrapText = $true

# Move the selection to match where the user ended up after typing all
# of this in (bottom of the new table).
$ws.Range("F88").Select()

Write-Host "done"
